# Apply edit: append rows 22-24 with new PGD-clipped results data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(20, 3, 0, 0.003, "Regular", "<function relu at 0x11c4be9d8>",
      0.9559000134468079, 0.05790000036358833, 0.04890000075101852,
      0.1738029420375824, 6.665348529815674, 0.05790000036358833,
      "logs/results_222.log", "weights/model_222.ckpt", "tb/222/non_robust",
      "(7.230026, 8.443772, 9.982576, 10.426454, 10.48337, 10.4655, 8.748777)",
      "(139.92761, 9.773078, 8.607171, 8.341426, 7.4237742, 7.5570107, 6.5512905, 8.353663)"),
    @(21, 3, 0, 0.003, "Regular", "<function relu at 0x11e21a9d8>",
      0.9312000274658203, 0.06909999996423721, 0.1206000000238419,
      0.2596416771411896, 7.456814289093018, 0.06909999996423721,
      "logs/results_232.log", "weights/model_232.ckpt", "tb/232/non_robust",
      "(7.3609123, 7.8955326, 9.737668, 9.91078, 10.763038, 10.914308, 9.4478855)",
      "(140.87561, 8.819877, 8.838883, 8.998733, 7.921186, 7.691656, 7.321085, 9.176294)"),
    @(22, 3, 0, 0.003, "Regular", "<function relu at 0x111a129d8>",
      0.9312999844551086, 0.03830000013113022, 0.009100000374019146,
      0.2439231872558594, 7.289804458618164, 0.03830000013113022,
      "logs/results_234.log", "weights/model_234.ckpt", "tb/234/non_robust",
      "(7.0172796, 7.6613173, 8.302304, 8.662146, 7.9462867, 9.3396435, 7.44935)",
      "(136.32106, 9.159889, 8.33736, 8.452992, 7.9244547, 7.2932935, 7.7801814, 9.626968)")
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

Write-Host "Done writing rows 22-24"
Write-Host "New UsedRange:" $ws.UsedRange.Address()
